$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value = 8630
$ws.Range("C6:C15").Value = 8293
$ws.Range("C16:C28").Value = 8201
$ws.Range("C29:C252").Value = 7569
